$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append a new job posting row (JD_020) right after the last existing data row.
$newRow = 21

$ws.Cells.Item($newRow, 1).Value = "JD_020"
$ws.Cells.Item($newRow, 2).Value = "Senior System Analyst"
$ws.Cells.Item($newRow, 3).Value = "We are seeking a Junior RPA Developer to design, develop, and support automation solutions.
Collaborate with teams to streamline business processes using RPA tools like UiPath or Automation Anywhere. Join Akkodis to grow your skills in a dynamic, tech-driven environment"
$ws.Cells.Item($newRow, 4).Value = 2
$ws.Cells.Item($newRow, 5).Value = 3

# Keep the new row's height consistent with the rest of the sheet (avoid Excel's
# automatic row-height recalculation triggered by the multi-line description text).
$ws.Rows.Item($newRow).AutoFit()
